# Generate Report for Handoff
# Update the "b9d5d801-102b-48d6-9175-1fbd789d7af6" rows across the
# Overview, zh-cn and de-de sheets to reflect that the file is now
# "Ready for handoff" with an updated handoff timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (row 3: b9d5d801-102b-48d6-9175-1fbd789d7af6.md) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-38-19 08:38:59"

# --- zh-cn sheet (row 3: b9d5d801-102b-48d6-9175-1fbd789d7af6.md) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-19 08:38:56"

# --- de-de sheet (row 3: b9d5d801-102b-48d6-9175-1fbd789d7af6.md) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-19 08:38:59"
